# Applies the "Update gh-pages to output generated at 456a3b4" change:
#  - Sheet "展览": refresh "想去人数" (F) counts for several events
#  - Sheet "演出": refresh one F value, and insert a new event row
#    ("杭州·英雄时代2024哈瓦西钢琴演奏会") before "杭州·苗阜王声 ..."
#  - Sheet "全部类型": refresh "想去人数" (F) counts for several events

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" F-column updates ---
$ws1.Range("F3").Value = 540
$ws1.Range("F5").Value = 502
$ws1.Range("F6").Value = 1138
$ws1.Range("F9").Value = 119
$ws1.Range("F10").Value = 120
$ws1.Range("F11").Value = 1172
$ws1.Range("F14").Value = 820
$ws1.Range("F15").Value = 839
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 70
$ws1.Range("F19").Value = 663
$ws1.Range("F20").Value = 709
$ws1.Range("F22").Value = 2534
$ws1.Range("F23").Value = 711
$ws1.Range("F24").Value = 77
$ws1.Range("F25").Value = 1974
$ws1.Range("F26").Value = 439
$ws1.Range("F27").Value = 2865
$ws1.Range("F28").Value = 530
$ws1.Range("F29").Value = 83
$ws1.Range("F30").Value = 698
$ws1.Range("F31").Value = 134
$ws1.Range("F32").Value = 111
$ws1.Range("F33").Value = 98
$ws1.Range("F34").Value = 995
$ws1.Range("F35").Value = 1717
$ws1.Range("F36").Value = 353
$ws1.Range("F38").Value = 541
$ws1.Range("F39").Value = 165
$ws1.Range("F40").Value = 124
$ws1.Range("F41").Value = 160
$ws1.Range("F42").Value = 24

# --- Sheet "演出" updates ---
$ws2.Range("F4").Value = 11

# Insert a brand-new row 13 (pushes the old rows 13 & 14 down to 14 & 15)
$ws2.Rows.Item(13).Insert()

# Copy column-A formatting (bold / border / centered, style used by every
# other row) from the row right below onto the freshly inserted, blank row.
$ws2.Range("A14").Copy()
$ws2.Range("A13").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Range("A13").Value = 12

# B13 ("2024-06-05") looks like a date, so force a text format first to
# keep it stored as plain text instead of being converted to a date
# serial number, then drop back to the Normal style so no extra
# number-format styling is left behind on the cell.
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2024-06-05"
$ws2.Range("B13").Style = "Normal"

$ws2.Range("C13").Value = "杭州·英雄时代2024哈瓦西钢琴演奏会"
$ws2.Range("D13").Value = "中国杭州北山路86号西湖岳湖景区 中国杭州西湖岳湖景区印象西湖"
$ws2.Range("E13").Value = "2024.06.05 20:00-06.05 21:30"
$ws2.Range("F13").Value = 0
$ws2.Range("G13").Value = 499
$ws2.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=83902"
$ws2.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202404/BFRFmKpT1712569969076.jpeg"

# The "序号" column (A) is a plain sequential index, not a formula, so a
# straight row-insert leaves the rows that got pushed down carrying their
# old index. Renumber them to stay sequential with the newly inserted row.
$ws2.Range("A14").Value = 13
$ws2.Range("A15").Value = 14

# --- Sheet "全部类型" F-column updates ---
$ws4.Range("F4").Value = 540
$ws4.Range("F6").Value = 502
$ws4.Range("F7").Value = 1138
$ws4.Range("F10").Value = 119
$ws4.Range("F11").Value = 120
$ws4.Range("F12").Value = 1173
$ws4.Range("F14").Value = 820
$ws4.Range("F15").Value = 839
$ws4.Range("F17").Value = 11
$ws4.Range("F18").Value = 58
$ws4.Range("F20").Value = 70
$ws4.Range("F21").Value = 709
$ws4.Range("F23").Value = 2534
$ws4.Range("F24").Value = 711
$ws4.Range("F25").Value = 77
$ws4.Range("F28").Value = 2865
$ws4.Range("F29").Value = 530
$ws4.Range("F34").Value = 83
$ws4.Range("F36").Value = 698
$ws4.Range("F37").Value = 134
$ws4.Range("F38").Value = 111
$ws4.Range("F39").Value = 98
$ws4.Range("F40").Value = 995
$ws4.Range("F41").Value = 1717
$ws4.Range("F43").Value = 353
$ws4.Range("F44").Value = 541
$ws4.Range("F45").Value = 165
$ws4.Range("F46").Value = 124
$ws4.Range("F47").Value = 160
$ws4.Range("F48").Value = 24
